$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Update the Version and Date values
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row before row 11 ("Description") to make room for "Jurisdiction",
# then copy the formatting from the row below so the inserted row keeps the
# same style as the rest of the table.
$meta.Rows.Item(11).Insert()
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New Jurisdiction row content. The value cell is an explicit empty string
# (not merely a blank cell), matching the other "no value" markers already
# used throughout this workbook (e.g. Elements!D2). A plain "" assignment
# collapses to a blank cell, so force text-empty via a quote-prefix, then
# re-apply the row's normal formatting to drop the quote-prefix style bit.
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = "'"
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
